# Edit: 
#  1) Swap the deck's colour theme: the slide master's theme colour scheme
#     changes from the "Integral" (Red Violet) palette to the default
#     "Office" palette.
#  2) Three tables (on the slides that hold a single data table each)
#     get their table style switched from the plain "Table_0" custom
#     style to the built-in Medium-style accent-1 table style.

$p = $ppt.ActivePresentation

# --- 1) Theme colour scheme -------------------------------------------------
$master = $p.Slides.Item(1).Master
$cs = $master.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0         # dk1      #000000
$cs.Item(2).RGB  = 16777215  # lt1      #FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      #44546A
$cs.Item(4).RGB  = 15132391  # lt2      #E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  #5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  #ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  #A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  #FFC000
$cs.Item(9).RGB  = 12874308  # accent5  #4472C4
$cs.Item(10).RGB = 4697456   # accent6  #70AD47
$cs.Item(11).RGB = 12673797  # hlink    #0563C1
$cs.Item(12).RGB = 7491477   # folHlink #954F72

# --- 2) Table styles ---------------------------------------------------------
$newStyleId = "{57060867-67A4-47DB-8B77-9F29A66E0EBA}"

foreach ($slideIndex in 14,15,16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
